# Apply the data update described by the commit:
#  - Header renamed from "Numbers" to "Test_Column"
#  - Numeric values (0,1,2) replaced with letters (a,b,c,d,e), extending
#    the column down to row 6
#  - Column A resized to fit the new content
#  - Selection moved to H16:H17
#  - Window position/size tweaked

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old range and write the new header + values
$ws.Range("A1").Value = "Test_Column"
$ws.Range("A2").Value = "a"
$ws.Range("A3").Value = "b"
$ws.Range("A4").Value = "c"
$ws.Range("A5").Value = "d"
$ws.Range("A6").Value = "e"

# Autofit column A to match the new "bestFit" width
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update the selection
$ws.Range("H16:H17").Select()
